$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header D1
$ws.Range("D1").Value = "ActionKeyword"

# Shift rows up: remove "Open the Browser" / "openBrowser" row, rename input_X -> inputX,
# and append "Close the browser" / "closeBrowser" as the new last row.
$ws.Range("C2").Value = "Navigate to website"
$ws.Range("D2").Value = "navigate"

$ws.Range("C3").Value = "Enter the Username in the Username field"
$ws.Range("D3").Value = "inputUsername"

$ws.Range("C4").Value = "Enter the Password in the Password field"
$ws.Range("D4").Value = "inputPassword"

$ws.Range("C5").Value = "Close the browser"
$ws.Range("D5").Value = "closeBrowser"

# Update selection to C5
$ws.Range("C5").Select()
